$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before the old column H ("Number of teachers" threshold block),
# shifting the existing H:M block to L:Q.
$ws.Range("H1:K1").EntireColumn.Insert()

# Populate the 4 new header cells -- the order below reproduces the shared-string
# insertion order seen in the target file (J, H, K, I).
$ws.Range("J1").Value = "Protection indicator value -- categorical variable"
$ws.Range("H1").Value = "Protection indicator value -- continuous / discrete numerical variable"
$ws.Range("K1").Value = "category for categorical variable --> if selected the child is in need"
$ws.Range("I1").Value = "threshold for numerical variable --> if above the child is in need"

# The 4 new header cells inherit the bold/centered/wrapped/bordered "Note" look from
# column G via the column insert; give them their own pale-yellow fill so they read as
# a visually distinct block of headers (mirrors the new gradient fill added upstream).
$headers = $ws.Range("H1:K1")
$headers.Interior.Pattern = -4105
$headers.Interior.Color = 15073279

# Column widths for the 4 new columns (closest values reachable through the
# ColumnWidth/character-width API).
$ws.Range("H:H").ColumnWidth = 21
$ws.Range("I:I").ColumnWidth = 21.67
$ws.Range("J:J").ColumnWidth = 22.67
$ws.Range("K:K").ColumnWidth = 24.33

# Row 1 grows taller to fit the new wrapped headers.
$ws.Range("1:1").RowHeight = 60

# Selection moves from a single stray cell to "select all".
$ws.Cells.Select()

Write-Output "applied protection-indicator columns"
